# Add a cumulative "n" column (F = D + E, total sample size) to the
# "Mata_et_al(2011)" sheet. Rows 10-12 and 15-16 already carry manually
# entered totals in column F that don't follow the D+E pattern, so they
# are intentionally left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mata_et_al(2011)")

$rows = @(2,3,4,5,6,7,8,9,13,14,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32)

foreach ($r in $rows) {
    $cell = $ws.Range("F$r")
    $cell.Formula = "=D$r+E$r"
    $cell.NumberFormat = "0"
}

$ws.Activate()
$ws.Range("F17:F32").Select()

$wb.Save()
